# "10Th - MB for single stock and added new group"
#
# The sheet tracks analyst price-target history by date, one column per
# business day, most-recent date on the left (column B) through oldest on
# the right. This edit adds two new snapshot dates (Jun_26, Jun_27) as
# three new leading columns, and appends two new rows for a "Benchmark"
# and "Evercore ISI" group/analyst.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns in front of the existing date columns (old B:E
# shifts right to E:H) to make room for the two new snapshot dates.
$ws.Range("B1:D1").EntireColumn.Insert()

# New header dates for the freshly inserted columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1:D1").Value = "Jun_26"

# Every analyst row (2-27) gets a rating of "UN" in the new columns, same
# as all of their other historical columns.
$ws.Range("B2:D27").Value = "UN"

# Add the new "Benchmark" group as row 28.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

# Add the new "Evercore ISI" analyst as row 29.
$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"
